$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.874436
$ws.Cells.Item(2, 8).Value = 5.623308
$ws.Cells.Item(2, 9).Value = 0.1442186763702422
$ws.Cells.Item(2, 10).Value = 0.1442186763702422
$ws.Cells.Item(2, 13).Value = 82.43338033333333
$ws.Cells.Item(2, 14).Value = 247.300141
$ws.Cells.Item(2, 15).Value = 0.3670006993429558
$ws.Cells.Item(2, 16).Value = 0.3670006993429557
$ws.Cells.Item(2, 17).Value = 154.516095698492
$ws.Cells.Item(2, 18).Value = 1390.644861286428
$ws.Cells.Item(2, 19).Value = 0.05292835508619429
$ws.Cells.Item(2, 20).Value = 0.05292835508619427
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.874436
$ws.Cells.Item(3, 8).Value = 5.623308
$ws.Cells.Item(3, 9).Value = 0.1442186763702422
$ws.Cells.Item(3, 10).Value = 0.1442186763702422
$ws.Cells.Item(3, 15).Value = 0.3956886215996139
$ws.Cells.Item(3, 16).Value = 0.3956886215996139
$ws.Cells.Item(3, 17).Value = 166.594398951692
$ws.Cells.Item(3, 18).Value = 1499.349590565228
$ws.Cells.Item(3, 19).Value = 0.05706568926186194
$ws.Cells.Item(3, 20).Value = 0.05706568926186192
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.874436
$ws.Cells.Item(4, 8).Value = 5.623308
$ws.Cells.Item(4, 9).Value = 0.1442186763702422
$ws.Cells.Item(4, 10).Value = 0.1442186763702422
$ws.Cells.Item(4, 13).Value = 42.93483766666667
$ws.Cells.Item(4, 14).Value = 128.804513
$ws.Cells.Item(4, 15).Value = 0.1911496942879982
$ws.Cells.Item(4, 16).Value = 0.1911496942879981
$ws.Cells.Item(4, 17).Value = 80.47860537655602
$ws.Cells.Item(4, 18).Value = 724.3074483890041
$ws.Cells.Item(4, 19).Value = 0.02756735589879154
$ws.Cells.Item(4, 20).Value = 0.02756735589879153
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.874436
$ws.Cells.Item(5, 8).Value = 5.623308
$ws.Cells.Item(5, 9).Value = 0.1442186763702422
$ws.Cells.Item(5, 10).Value = 0.1442186763702422
$ws.Cells.Item(5, 13).Value = 10.368389
$ws.Cells.Item(5, 14).Value = 31.105167
$ws.Cells.Item(5, 15).Value = 0.04616098476943217
$ws.Cells.Item(5, 16).Value = 0.04616098476943217
$ws.Cells.Item(5, 17).Value = 19.434881603604
$ws.Cells.Item(5, 18).Value = 174.913934432436
$ws.Cells.Item(5, 19).Value = 0.006657276123394417
$ws.Cells.Item(5, 20).Value = 0.006657276123394415
$ws.Cells.Item(6, 9).Value = 0.2460517715407892
$ws.Cells.Item(6, 10).Value = 0.2460517715407892
$ws.Cells.Item(6, 13).Value = 82.43338033333333
$ws.Cells.Item(6, 14).Value = 247.300141
$ws.Cells.Item(6, 15).Value = 0.3670006993429558
$ws.Cells.Item(6, 16).Value = 0.3670006993429557
$ws.Cells.Item(6, 17).Value = 263.620219205013
$ws.Cells.Item(6, 18).Value = 2372.581972845117
$ws.Cells.Item(6, 19).Value = 0.09030117223004282
$ws.Cells.Item(6, 20).Value = 0.09030117223004279
$ws.Cells.Item(7, 9).Value = 0.2460517715407892
$ws.Cells.Item(7, 10).Value = 0.2460517715407892
$ws.Cells.Item(7, 15).Value = 0.3956886215996139
$ws.Cells.Item(7, 16).Value = 0.3956886215996139
$ws.Cells.Item(7, 19).Value = 0.097359886323118
$ws.Cells.Item(7, 20).Value = 0.09735988632311797
$ws.Cells.Item(8, 9).Value = 0.2460517715407892
$ws.Cells.Item(8, 10).Value = 0.2460517715407892
$ws.Cells.Item(8, 13).Value = 42.93483766666667
$ws.Cells.Item(8, 14).Value = 128.804513
$ws.Cells.Item(8, 15).Value = 0.1911496942879982
$ws.Cells.Item(8, 16).Value = 0.1911496942879981
$ws.Cells.Item(8, 17).Value = 137.304709226409
$ws.Cells.Item(8, 18).Value = 1235.742383037681
$ws.Cells.Item(8, 19).Value = 0.04703272090904222
$ws.Cells.Item(8, 20).Value = 0.04703272090904221
$ws.Cells.Item(9, 9).Value = 0.2460517715407892
$ws.Cells.Item(9, 10).Value = 0.2460517715407892
$ws.Cells.Item(9, 13).Value = 10.368389
$ws.Cells.Item(9, 14).Value = 31.105167
$ws.Cells.Item(9, 15).Value = 0.04616098476943217
$ws.Cells.Item(9, 16).Value = 0.04616098476943217
$ws.Cells.Item(9, 17).Value = 33.157890285831
$ws.Cells.Item(9, 18).Value = 298.421012572479
$ws.Cells.Item(9, 19).Value = 0.01135799207858617
$ws.Cells.Item(9, 20).Value = 0.01135799207858617
$ws.Cells.Item(10, 7).Value = 6.825289333333334
$ws.Cells.Item(10, 8).Value = 20.475868
$ws.Cells.Item(10, 9).Value = 0.5251361975000832
$ws.Cells.Item(10, 10).Value = 0.5251361975000832
$ws.Cells.Item(10, 13).Value = 82.43338033333333
$ws.Cells.Item(10, 14).Value = 247.300141
$ws.Cells.Item(10, 15).Value = 0.3670006993429558
$ws.Cells.Item(10, 16).Value = 0.3670006993429557
$ws.Cells.Item(10, 17).Value = 562.6316714997098
$ws.Cells.Item(10, 18).Value = 5063.685043497389
$ws.Cells.Item(10, 19).Value = 0.1927253517328311
$ws.Cells.Item(10, 20).Value = 0.1927253517328311
$ws.Cells.Item(11, 7).Value = 6.825289333333334
$ws.Cells.Item(11, 8).Value = 20.475868
$ws.Cells.Item(11, 9).Value = 0.5251361975000832
$ws.Cells.Item(11, 10).Value = 0.5251361975000832
$ws.Cells.Item(11, 15).Value = 0.3956886215996139
$ws.Cells.Item(11, 16).Value = 0.3956886215996139
$ws.Cells.Item(11, 17).Value = 606.6117883769098
$ws.Cells.Item(11, 18).Value = 5459.506095392188
$ws.Cells.Item(11, 19).Value = 0.2077904181408706
$ws.Cells.Item(11, 20).Value = 0.2077904181408705
$ws.Cells.Item(12, 7).Value = 6.825289333333334
$ws.Cells.Item(12, 8).Value = 20.475868
$ws.Cells.Item(12, 9).Value = 0.5251361975000832
$ws.Cells.Item(12, 10).Value = 0.5251361975000832
$ws.Cells.Item(12, 13).Value = 42.93483766666667
$ws.Cells.Item(12, 14).Value = 128.804513
$ws.Cells.Item(12, 15).Value = 0.1911496942879982
$ws.Cells.Item(12, 16).Value = 0.1911496942879981
$ws.Cells.Item(12, 17).Value = 293.0426895546983
$ws.Cells.Item(12, 18).Value = 2637.384205992284
$ws.Cells.Item(12, 19).Value = 0.1003796236117027
$ws.Cells.Item(12, 20).Value = 0.1003796236117027
$ws.Cells.Item(13, 7).Value = 6.825289333333334
$ws.Cells.Item(13, 8).Value = 20.475868
$ws.Cells.Item(13, 9).Value = 0.5251361975000832
$ws.Cells.Item(13, 10).Value = 0.5251361975000832
$ws.Cells.Item(13, 13).Value = 10.368389
$ws.Cells.Item(13, 14).Value = 31.105167
$ws.Cells.Item(13, 15).Value = 0.04616098476943217
$ws.Cells.Item(13, 16).Value = 0.04616098476943217
$ws.Cells.Item(13, 17).Value = 70.76725484555068
$ws.Cells.Item(13, 18).Value = 636.9052936099561
$ws.Cells.Item(13, 19).Value = 0.02424080401467887
$ws.Cells.Item(13, 20).Value = 0.02424080401467886
$ws.Cells.Item(14, 7).Value = 1.099475
$ws.Cells.Item(14, 8).Value = 3.298425
$ws.Cells.Item(14, 9).Value = 0.08459335458888541
$ws.Cells.Item(14, 10).Value = 0.08459335458888539
$ws.Cells.Item(14, 13).Value = 82.43338033333333
$ws.Cells.Item(14, 14).Value = 247.300141
$ws.Cells.Item(14, 15).Value = 0.3670006993429558
$ws.Cells.Item(14, 16).Value = 0.3670006993429557
$ws.Cells.Item(14, 17).Value = 90.63344084199166
$ws.Cells.Item(14, 18).Value = 815.700967577925
$ws.Cells.Item(14, 19).Value = 0.03104582029388758
$ws.Cells.Item(14, 20).Value = 0.03104582029388757
$ws.Cells.Item(15, 7).Value = 1.099475
$ws.Cells.Item(15, 8).Value = 3.298425
$ws.Cells.Item(15, 9).Value = 0.08459335458888541
$ws.Cells.Item(15, 10).Value = 0.08459335458888539
$ws.Cells.Item(15, 15).Value = 0.3956886215996139
$ws.Cells.Item(15, 16).Value = 0.3956886215996139
$ws.Cells.Item(15, 17).Value = 97.71812789949166
$ws.Cells.Item(15, 18).Value = 879.4631510954249
$ws.Cells.Item(15, 19).Value = 0.03347262787376344
$ws.Cells.Item(15, 20).Value = 0.03347262787376343
$ws.Cells.Item(16, 7).Value = 1.099475
$ws.Cells.Item(16, 8).Value = 3.298425
$ws.Cells.Item(16, 9).Value = 0.08459335458888541
$ws.Cells.Item(16, 10).Value = 0.08459335458888539
$ws.Cells.Item(16, 13).Value = 42.93483766666667
$ws.Cells.Item(16, 14).Value = 128.804513
$ws.Cells.Item(16, 15).Value = 0.1911496942879982
$ws.Cells.Item(16, 16).Value = 0.1911496942879981
$ws.Cells.Item(16, 17).Value = 47.20578064355834
$ws.Cells.Item(16, 18).Value = 424.852025792025
$ws.Cells.Item(16, 19).Value = 0.01616999386846167
$ws.Cells.Item(16, 20).Value = 0.01616999386846167
$ws.Cells.Item(17, 7).Value = 1.099475
$ws.Cells.Item(17, 8).Value = 3.298425
$ws.Cells.Item(17, 9).Value = 0.08459335458888541
$ws.Cells.Item(17, 10).Value = 0.08459335458888539
$ws.Cells.Item(17, 13).Value = 10.368389
$ws.Cells.Item(17, 14).Value = 31.105167
$ws.Cells.Item(17, 15).Value = 0.04616098476943217
$ws.Cells.Item(17, 16).Value = 0.04616098476943217
$ws.Cells.Item(17, 17).Value = 11.399784495775
$ws.Cells.Item(17, 18).Value = 102.598060461975
$ws.Cells.Item(17, 19).Value = 0.003904912552772715
$ws.Cells.Item(17, 20).Value = 0.003904912552772713
